$d = $word.ActiveDocument

# Helper: split off the trailing sub-range [splitStart, rangeEnd) into its own
# run by toggling Bold on/off (forces a run boundary without altering the
# visible formatting of the text).
function Split-Run($splitStart, $rangeEnd) {
    $wr = $d.Range($splitStart, $rangeEnd)
    $wr.Font.Bold = 1
    $wr.Font.Bold = 0
}

# Helper: given the full text of a sentence/phrase (already present verbatim
# in the document), locate it and force a run boundary right before its last
# word - this is where Word's grammar checker would have wrapped the word in
# <w:proofErr w:type="gramStart"/>...<w:proofErr w:type="gramEnd"/>.
function Split-LastWord($sentenceText) {
    $full = $d.Content
    $found = $full.Find.Execute($sentenceText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        return
    }
    $sentEnd = $full.End
    $words = $full.Words
    $lastWord = $words.Item($words.Count - 1)
    Split-Run $lastWord.Start $sentEnd
}

# --- 1. LinkedIn URL: split "pranavkhismatrao" into its own run (spell-check flag) ---
$full = $d.Content
$found = $full.Find.Execute("pranavkhismatrao", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $wordStart = $full.Start
    $wordEnd = $full.End
    Split-Run $wordStart $wordEnd
    $full2 = $d.Content
    $full2.Find.Execute("pranavkhismatrao/", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    Split-Run $wordEnd $full2.End
}

# --- 2. Remove ", Proficient in Math and Linear Algebra." ---
$d.Content.Find.Execute(", Proficient in Math and Linear Algebra.", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# --- 3. Frameworks line: insert ", Angular, React, Vue" in place of the first two tabs after ", .Net MVC" ---
$full = $d.Content
$full.Find.Execute(", .Net MVC", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insertPoint = $full.End
$tabRange = $d.Range($insertPoint, $insertPoint + 2)
$tabRange.Text = ", Angular, React, Vue"

# --- 4. "...and MVC framework" -> split off "framework" ---
Split-LastWord "Accomplished 6 months of training held by TCS with course content including C# language, and MVC framework"

# --- 5. "...utilizing agile methodology" -> split off "methodology" ---
Split-LastWord "Programmed 3 comprehensive solutions with a competent team of 3 leveraging RPA technology along with testing 20 deploy Cases utilizing agile methodology"

# --- 6. " using visual Studio" -> split off "Studio" ---
Split-LastWord " using visual Studio"

# --- 7. "...and double efficiency" -> split off "efficiency" ---
Split-LastWord "Created a web portal with a team of three to create, edit, read, and track FAQ tickets to cut down on lag time and double efficiency"

# --- 8. "with a team of five" -> split off "five" ---
Split-LastWord "with a team of five"

# --- 9. "...and during UAT" -> split off "UAT" ---
Split-LastWord "Utilizing critical thinking abilities to identify the source of problems encountered both during the development phase and during UAT"

# --- 10. "database" run gets wrapped (already its own run; nothing textual to change) ---

# --- 11. "...build the frontend" -> split off "frontend" ---
Split-LastWord " to call APIs on the backend while HTML, CSS, and Bootstrap were used to build the frontend"

# --- 12. "built" run gets wrapped (already its own run; nothing textual to change) ---

# --- 13. "Aerospike, API call handling" -> split off "handling" ---
Split-LastWord "Aerospike, API call handling"
